$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "36.305.64"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -2.71%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.972.31"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -3.33%  "

$ws.Range("E4").Value = "  -0.06%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "245.20"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.30%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.623"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.83%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "58.27"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -10.94%  "

$ws.Range("E8").Value = "  -0.03%  "

$ws.Range("E9").Value = "  -6.44%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "56.18"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -5.42%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0866"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +9.27%  "

$ws.Range("E12").Value = "  -0.14%  "

$ws.Range("E13").Value = "  -6.23%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.263.27"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.23%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "21.70"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -6.37%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "13.65"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -7.36%  "

$ws.Range("E17").Value = "  -4.89%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.985.01"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.77%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "36.196.39"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.65%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0897"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.75%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "70.16"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.84%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.22"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -4.49%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "233.66"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.91%  "

$ws.Range("E24").Value = "  +0.13%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.48"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -4.05%  "

$ws.Range("E26").Value = "  -2.45%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.74"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.52%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "165.87"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.04%  "

$ws.Range("E29").Value = "  -0.16%  "

$ws.Range("E30").Value = "  +1.59%  "

$ws.Range("E31").Value = "  -2.19%  "

$ws.Range("E32").Value = "  -1.40%  "

$ws.Range("E33").Value = "  -5.36%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0643"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.04%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.38"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -5.21%  "

$ws.Range("E37").Value = "  -5.22%  "

$ws.Range("E38").Value = "  -1.85%  "

$ws.Range("E39").Value = "  -6.72%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.91"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.64%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0958"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -5.60%  "

$ws.Range("E42").Value = "  -6.92%  "

$ws.Range("E43").Value = "  -5.00%  "

$ws.Range("E44").Value = "  -2.99%  "

$ws.Range("E45").Value = "  -7.15%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "15.95"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -7.93%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "90.43"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -4.95%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.358.27"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.52%  "

$ws.Range("E49").Value = "  -5.88%  "

$ws.Range("E50").Value = "  -2.87%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "45.15"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.88%  "

